$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "'PROC-2023-0001"
$ws.Range("I2").Style = "Normal"
$ws.Range("J2").Value = "'2023-01-19"
$ws.Range("J2").Style = "Normal"
$ws.Range("L2").Value = "'2023-06-01"
$ws.Range("L2").Style = "Normal"
$ws.Range("M2").Value = "'Sim"
$ws.Range("M2").Style = "Normal"
$ws.Range("N2").Value = "'2025-06-26"
$ws.Range("N2").Style = "Normal"
$ws.Range("I3").Value = "'PROC-2023-0002"
$ws.Range("I3").Style = "Normal"
$ws.Range("J3").Value = "'2023-11-20"
$ws.Range("J3").Style = "Normal"
$ws.Range("K3").Value = "'Não"
$ws.Range("K3").Style = "Normal"
$ws.Range("M3").Value = "'Não"
$ws.Range("M3").Style = "Normal"
$ws.Range("I4").Value = "'PROC-2024-0003"
$ws.Range("I4").Style = "Normal"
$ws.Range("J4").Value = "'2024-10-05"
$ws.Range("J4").Style = "Normal"
$ws.Range("K4").Value = "'Não"
$ws.Range("K4").Style = "Normal"
$ws.Range("J5").Value = "'2023-06-24"
$ws.Range("J5").Style = "Normal"
$ws.Range("K5").Value = "'Sim"
$ws.Range("K5").Style = "Normal"
$ws.Range("L5").Value = "'2025-10-16"
$ws.Range("L5").Style = "Normal"
$ws.Range("J6").Value = "'2024-09-19"
$ws.Range("J6").Style = "Normal"
$ws.Range("L6").Value = "'2025-11-06"
$ws.Range("L6").Style = "Normal"
$ws.Range("N6").Value = "'2025-11-10"
$ws.Range("N6").Style = "Normal"
$ws.Range("J7").Value = "'2025-02-28"
$ws.Range("J7").Style = "Normal"
$ws.Range("K7").Value = "'Não"
$ws.Range("K7").Style = "Normal"
$ws.Range("M7").Value = "'Não"
$ws.Range("M7").Style = "Normal"
$ws.Range("I8").Value = "'PROC-2023-0007"
$ws.Range("I8").Style = "Normal"
$ws.Range("J8").Value = "'2023-06-15"
$ws.Range("J8").Style = "Normal"
$ws.Range("K8").Value = "'Sim"
$ws.Range("K8").Style = "Normal"
$ws.Range("L8").Value = "'2023-06-27"
$ws.Range("L8").Style = "Normal"
$ws.Range("M8").Value = "'Sim"
$ws.Range("M8").Style = "Normal"
$ws.Range("N8").Value = "'2025-01-05"
$ws.Range("N8").Style = "Normal"
$ws.Range("J9").Value = "'2023-09-15"
$ws.Range("J9").Style = "Normal"
$ws.Range("L9").Value = "'2025-04-09"
$ws.Range("L9").Style = "Normal"
$ws.Range("M9").Value = "'Sim"
$ws.Range("M9").Style = "Normal"
$ws.Range("N9").Value = "'2025-05-14"
$ws.Range("N9").Style = "Normal"
$ws.Range("I10").Value = "'PROC-2023-0009"
$ws.Range("I10").Style = "Normal"
$ws.Range("J10").Value = "'2023-04-09"
$ws.Range("J10").Style = "Normal"
$ws.Range("J11").Value = "'2023-05-31"
$ws.Range("J11").Style = "Normal"
$ws.Range("L11").Value = "'2025-01-25"
$ws.Range("L11").Style = "Normal"
$ws.Range("I12").Value = "'PROC-2024-0011"
$ws.Range("I12").Style = "Normal"
$ws.Range("J12").Value = "'2025-09-18"
$ws.Range("J12").Style = "Normal"
$ws.Range("J13").Value = "'2023-12-14"
$ws.Range("J13").Style = "Normal"
$ws.Range("L13").Value = "'2023-11-25"
$ws.Range("L13").Style = "Normal"
$ws.Range("I14").Value = "'PROC-2025-0013"
$ws.Range("I14").Style = "Normal"
$ws.Range("J14").Value = "'2024-10-13"
$ws.Range("J14").Style = "Normal"
$ws.Range("K14").Value = "'Sim"
$ws.Range("K14").Style = "Normal"
$ws.Range("L14").Value = "'2025-09-05"
$ws.Range("L14").Style = "Normal"
$ws.Range("I15").Value = "'PROC-2023-0014"
$ws.Range("I15").Style = "Normal"
$ws.Range("J15").Value = "'2024-11-18"
$ws.Range("J15").Style = "Normal"
$ws.Range("K15").Value = "'Sim"
$ws.Range("K15").Style = "Normal"
$ws.Range("L15").Value = "'2023-01-21"
$ws.Range("L15").Style = "Normal"
$ws.Range("M15").Value = "'Sim"
$ws.Range("M15").Style = "Normal"
$ws.Range("N15").Value = "'2024-01-13"
$ws.Range("N15").Style = "Normal"
$ws.Range("I16").Value = "'PROC-2025-0015"
$ws.Range("I16").Style = "Normal"
$ws.Range("J16").Value = "'2024-07-31"
$ws.Range("J16").Style = "Normal"
$ws.Range("K16").Value = "'Sim"
$ws.Range("K16").Style = "Normal"
$ws.Range("L16").Value = "'2025-12-04"
$ws.Range("L16").Style = "Normal"
$ws.Range("M16").Value = "'Sim"
$ws.Range("M16").Style = "Normal"
$ws.Range("N16").Value = "'2023-02-28"
$ws.Range("N16").Style = "Normal"
$ws.Range("I17").Value = "'PROC-2025-0016"
$ws.Range("I17").Style = "Normal"
$ws.Range("J17").Value = "'2023-03-23"
$ws.Range("J17").Style = "Normal"
$ws.Range("L17").Value = "'2024-09-12"
$ws.Range("L17").Style = "Normal"
$ws.Range("M17").Value = "'Não"
$ws.Range("M17").Style = "Normal"
$ws.Range("J18").Value = "'2024-03-24"
$ws.Range("J18").Style = "Normal"
$ws.Range("I19").Value = "'PROC-2025-0018"
$ws.Range("I19").Style = "Normal"
$ws.Range("J19").Value = "'2024-05-10"
$ws.Range("J19").Style = "Normal"
$ws.Range("K19").Value = "'Não"
$ws.Range("K19").Style = "Normal"
$ws.Range("J20").Value = "'2025-04-11"
$ws.Range("J20").Style = "Normal"
$ws.Range("L20").Value = "'2024-04-17"
$ws.Range("L20").Style = "Normal"
$ws.Range("N20").Value = "'2024-08-28"
$ws.Range("N20").Style = "Normal"
$ws.Range("J21").Value = "'2025-09-01"
$ws.Range("J21").Style = "Normal"
$ws.Range("L21").Value = "'2023-11-08"
$ws.Range("L21").Style = "Normal"
$ws.Range("N21").Value = "'2023-04-02"
$ws.Range("N21").Style = "Normal"
$ws.Range("I22").Value = "'PROC-2023-0021"
$ws.Range("I22").Style = "Normal"
$ws.Range("J22").Value = "'2025-05-28"
$ws.Range("J22").Style = "Normal"
$ws.Range("K22").Value = "'Sim"
$ws.Range("K22").Style = "Normal"
$ws.Range("L22").Value = "'2023-02-06"
$ws.Range("L22").Style = "Normal"
$ws.Range("M22").Value = "'Sim"
$ws.Range("M22").Style = "Normal"
$ws.Range("N22").Value = "'2025-04-29"
$ws.Range("N22").Style = "Normal"
$ws.Range("I23").Value = "'PROC-2025-0022"
$ws.Range("I23").Style = "Normal"
$ws.Range("J23").Value = "'2025-02-14"
$ws.Range("J23").Style = "Normal"
$ws.Range("I24").Value = "'PROC-2025-0023"
$ws.Range("I24").Style = "Normal"
$ws.Range("J24").Value = "'2024-01-08"
$ws.Range("J24").Style = "Normal"
$ws.Range("K24").Value = "'Não"
$ws.Range("K24").Style = "Normal"
$ws.Range("M24").Value = "'Não"
$ws.Range("M24").Style = "Normal"
$ws.Range("I25").Value = "'PROC-2025-0024"
$ws.Range("I25").Style = "Normal"
$ws.Range("J25").Value = "'2025-04-17"
$ws.Range("J25").Style = "Normal"
$ws.Range("K25").Value = "'Sim"
$ws.Range("K25").Style = "Normal"
$ws.Range("L25").Value = "'2023-05-18"
$ws.Range("L25").Style = "Normal"
$ws.Range("I26").Value = "'PROC-2024-0025"
$ws.Range("I26").Style = "Normal"
$ws.Range("J26").Value = "'2023-03-25"
$ws.Range("J26").Style = "Normal"
$ws.Range("L26").Value = "'2025-04-27"
$ws.Range("L26").Style = "Normal"
$ws.Range("N26").Value = "'2024-12-14"
$ws.Range("N26").Style = "Normal"
$ws.Range("I27").Value = "'PROC-2024-0026"
$ws.Range("I27").Style = "Normal"
$ws.Range("J27").Value = "'2023-06-09"
$ws.Range("J27").Style = "Normal"
$ws.Range("L27").Value = "'2025-01-06"
$ws.Range("L27").Style = "Normal"
$ws.Range("N27").Value = "'2024-02-17"
$ws.Range("N27").Style = "Normal"
$ws.Range("I28").Value = "'PROC-2024-0027"
$ws.Range("I28").Style = "Normal"
$ws.Range("J28").Value = "'2025-10-24"
$ws.Range("J28").Style = "Normal"
$ws.Range("K28").Value = "'Não"
$ws.Range("K28").Style = "Normal"
$ws.Range("M28").Value = "'Não"
$ws.Range("M28").Style = "Normal"
$ws.Range("J29").Value = "'2024-07-22"
$ws.Range("J29").Style = "Normal"
$ws.Range("I30").Value = "'PROC-2023-0029"
$ws.Range("I30").Style = "Normal"
$ws.Range("J30").Value = "'2023-01-06"
$ws.Range("J30").Style = "Normal"
$ws.Range("K30").Value = "'Sim"
$ws.Range("K30").Style = "Normal"
$ws.Range("L30").Value = "'2024-10-05"
$ws.Range("L30").Style = "Normal"
$ws.Range("M30").Value = "'Sim"
$ws.Range("M30").Style = "Normal"
$ws.Range("N30").Value = "'2024-09-22"
$ws.Range("N30").Style = "Normal"
$ws.Range("I31").Value = "'PROC-2024-0030"
$ws.Range("I31").Style = "Normal"
$ws.Range("J31").Value = "'2024-10-07"
$ws.Range("J31").Style = "Normal"
$ws.Range("L31").Value = "'2024-05-24"
$ws.Range("L31").Style = "Normal"
$ws.Range("M31").Value = "'Sim"
$ws.Range("M31").Style = "Normal"
$ws.Range("N31").Value = "'2025-03-14"
$ws.Range("N31").Style = "Normal"
$ws.Range("I32").Value = "'PROC-2025-0031"
$ws.Range("I32").Style = "Normal"
$ws.Range("J32").Value = "'2025-02-28"
$ws.Range("J32").Style = "Normal"
$ws.Range("L32").Value = "'2024-07-15"
$ws.Range("L32").Style = "Normal"
$ws.Range("I33").Value = "'PROC-2023-0032"
$ws.Range("I33").Style = "Normal"
$ws.Range("J33").Value = "'2023-08-10"
$ws.Range("J33").Style = "Normal"
$ws.Range("K33").Value = "'Sim"
$ws.Range("K33").Style = "Normal"
$ws.Range("L33").Value = "'2024-02-21"
$ws.Range("L33").Style = "Normal"
$ws.Range("M33").Value = "'Sim"
$ws.Range("M33").Style = "Normal"
$ws.Range("N33").Value = "'2023-07-24"
$ws.Range("N33").Style = "Normal"
$ws.Range("I34").Value = "'PROC-2025-0033"
$ws.Range("I34").Style = "Normal"
$ws.Range("J34").Value = "'2023-11-25"
$ws.Range("J34").Style = "Normal"
$ws.Range("K34").Value = "'Não"
$ws.Range("K34").Style = "Normal"
$ws.Range("M34").Value = "'Não"
$ws.Range("M34").Style = "Normal"
$ws.Range("I35").Value = "'PROC-2025-0034"
$ws.Range("I35").Style = "Normal"
$ws.Range("J35").Value = "'2024-07-13"
$ws.Range("J35").Style = "Normal"
$ws.Range("L35").Value = "'2023-04-16"
$ws.Range("L35").Style = "Normal"
$ws.Range("M35").Value = "'Sim"
$ws.Range("M35").Style = "Normal"
$ws.Range("N35").Value = "'2025-10-07"
$ws.Range("N35").Style = "Normal"
$ws.Range("J36").Value = "'2024-05-28"
$ws.Range("J36").Style = "Normal"
$ws.Range("L36").Value = "'2023-01-03"
$ws.Range("L36").Style = "Normal"
$ws.Range("M36").Value = "'Não"
$ws.Range("M36").Style = "Normal"
$ws.Range("I37").Value = "'PROC-2024-0036"
$ws.Range("I37").Style = "Normal"
$ws.Range("J37").Value = "'2024-02-01"
$ws.Range("J37").Style = "Normal"
$ws.Range("L37").Value = "'2025-11-13"
$ws.Range("L37").Style = "Normal"
$ws.Range("N37").Value = "'2023-04-24"
$ws.Range("N37").Style = "Normal"
$ws.Range("I38").Value = "'PROC-2024-0037"
$ws.Range("I38").Style = "Normal"
$ws.Range("J38").Value = "'2025-08-07"
$ws.Range("J38").Style = "Normal"
$ws.Range("K38").Value = "'Sim"
$ws.Range("K38").Style = "Normal"
$ws.Range("L38").Value = "'2025-09-15"
$ws.Range("L38").Style = "Normal"
$ws.Range("M38").Value = "'Sim"
$ws.Range("M38").Style = "Normal"
$ws.Range("N38").Value = "'2023-11-04"
$ws.Range("N38").Style = "Normal"
$ws.Range("J39").Value = "'2024-01-21"
$ws.Range("J39").Style = "Normal"
$ws.Range("L39").Value = "'2024-04-20"
$ws.Range("L39").Style = "Normal"
$ws.Range("N39").Value = "'2023-01-27"
$ws.Range("N39").Style = "Normal"
$ws.Range("I40").Value = "'PROC-2023-0039"
$ws.Range("I40").Style = "Normal"
$ws.Range("J40").Value = "'2024-09-17"
$ws.Range("J40").Style = "Normal"
$ws.Range("L40").Value = "'2025-06-20"
$ws.Range("L40").Style = "Normal"
$ws.Range("M40").Value = "'Não"
$ws.Range("M40").Style = "Normal"
$ws.Range("I41").Value = "'PROC-2025-0040"
$ws.Range("I41").Style = "Normal"
$ws.Range("J41").Value = "'2025-03-21"
$ws.Range("J41").Style = "Normal"
$ws.Range("L41").Value = "'2024-01-17"
$ws.Range("L41").Style = "Normal"
$ws.Range("M41").Value = "'Sim"
$ws.Range("M41").Style = "Normal"
$ws.Range("N41").Value = "'2023-03-24"
$ws.Range("N41").Style = "Normal"
$ws.Range("I42").Value = "'PROC-2023-0041"
$ws.Range("I42").Style = "Normal"
$ws.Range("J42").Value = "'2024-09-28"
$ws.Range("J42").Style = "Normal"
$ws.Range("L42").Value = "'2024-05-30"
$ws.Range("L42").Style = "Normal"
$ws.Range("M42").Value = "'Sim"
$ws.Range("M42").Style = "Normal"
$ws.Range("N42").Value = "'2024-03-18"
$ws.Range("N42").Style = "Normal"
$ws.Range("I43").Value = "'PROC-2023-0042"
$ws.Range("I43").Style = "Normal"
$ws.Range("J43").Value = "'2025-12-21"
$ws.Range("J43").Style = "Normal"
$ws.Range("J44").Value = "'2023-11-08"
$ws.Range("J44").Style = "Normal"
$ws.Range("K44").Value = "'Não"
$ws.Range("K44").Style = "Normal"
$ws.Range("M44").Value = "'Não"
$ws.Range("M44").Style = "Normal"
$ws.Range("I45").Value = "'PROC-2023-0044"
$ws.Range("I45").Style = "Normal"
$ws.Range("J45").Value = "'2024-09-04"
$ws.Range("J45").Style = "Normal"
$ws.Range("L45").Value = "'2023-04-22"
$ws.Range("L45").Style = "Normal"
$ws.Range("N45").Value = "'2025-09-04"
$ws.Range("N45").Style = "Normal"
$ws.Range("J46").Value = "'2024-04-15"
$ws.Range("J46").Style = "Normal"
$ws.Range("L46").Value = "'2023-05-11"
$ws.Range("L46").Style = "Normal"
$ws.Range("M46").Value = "'Não"
$ws.Range("M46").Style = "Normal"
$ws.Range("I47").Value = "'PROC-2023-0046"
$ws.Range("I47").Style = "Normal"
$ws.Range("J47").Value = "'2024-06-21"
$ws.Range("J47").Style = "Normal"
$ws.Range("K47").Value = "'Não"
$ws.Range("K47").Style = "Normal"
$ws.Range("M47").Value = "'Não"
$ws.Range("M47").Style = "Normal"
$ws.Range("J48").Value = "'2024-08-19"
$ws.Range("J48").Style = "Normal"
$ws.Range("K48").Value = "'Não"
$ws.Range("K48").Style = "Normal"
$ws.Range("I49").Value = "'PROC-2025-0048"
$ws.Range("I49").Style = "Normal"
$ws.Range("J49").Value = "'2023-12-17"
$ws.Range("J49").Style = "Normal"
$ws.Range("K49").Value = "'Não"
$ws.Range("K49").Style = "Normal"
$ws.Range("J50").Value = "'2024-01-24"
$ws.Range("J50").Style = "Normal"
$ws.Range("I51").Value = "'PROC-2025-0050"
$ws.Range("I51").Style = "Normal"
$ws.Range("J51").Value = "'2024-03-08"
$ws.Range("J51").Style = "Normal"
$ws.Range("L51").Value = "'2023-06-18"
$ws.Range("L51").Style = "Normal"
$ws.Range("N51").Value = "'2023-11-03"
$ws.Range("N51").Style = "Normal"

$ws.Range("L3").ClearContents()
$ws.Range("N3").ClearContents()
$ws.Range("L4").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("N17").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("L24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("L28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("L34").ClearContents()
$ws.Range("N34").ClearContents()
$ws.Range("N36").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("L44").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("N46").ClearContents()
$ws.Range("L47").ClearContents()
$ws.Range("N47").ClearContents()
$ws.Range("L48").ClearContents()
$ws.Range("L49").ClearContents()
